# Update Issue 1 (Deposit an Product)
# - Rename sheet "Report Growth" -> "Report Summary"
# - Relabel header cells with "(Rp)" currency suffix
# - Move active selection from G6 to D1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Report Summary"

# Update header row labels (column A "Location" stays the same)
$ws.Range("B1").Value = "Return (Rp)"
$ws.Range("C1").Value = "Used (Rp)"
$ws.Range("D1").Value = "Remaining (Rp)"

# Update the selected/active cell shown when the file is reopened
$ws.Range("D1").Select()
